# Helios Matrix - balancing update
# - Updated Colliders
# - Updated Decks
# - Balancing update

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Balancing value updates (column P = "Water reservoir" collider costs, etc.)
$ws.Range("P6").Value = 1
$ws.Range("P7").Value = 1
$ws.Range("P8").Value = 1
$ws.Range("P9").Value = 1
$ws.Range("P10").Value = 10
$ws.Range("N11").Value = 5
$ws.Range("P16").Value = 2
$ws.Range("P17").Value = 1
$ws.Range("Q17").Value = 7

# Row 17 height tweak (explicit custom height instead of default 16.5)
$ws.Rows.Item(17).RowHeight = 15.75

# Move the active selection to P24
$ws.Activate()
$ws.Range("P24").Select()
